# Update the crypto price/volume table (columns D = Price, E = Volume(1h))
# with refreshed figures, matching the "Updated cryptos list" GitHub Actions
# commit. Numeric-looking price strings are written via NumberFormat "@"
# (Text) so Excel keeps them as literal text instead of coercing them into
# floating point numbers (which would both change the cell type and lose
# precision/formatting, e.g. "594.59" -> 594.59000000000003). The style is
# reset back to "Normal" afterwards so no lingering direct formatting is
# left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.092.22'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '3.516.88'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +6.84%  '
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').Value = '4.121.81'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.134'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.78'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.08%  '
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').Value = '67.072.45'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').Value = '3.524.34'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '394.90'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.23'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('E25').Value = '  -4.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('E29').Value = '  -2.22%  '
$ws.Range('E30').Value = '  -1.33%  '
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.53%  '
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '163.65'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.895'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.90'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.66%  '
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0742'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '26.47'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '27.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.62'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.07%  '
$ws.Range('D44').Value = '2.797.20'
$ws.Range('E44').Value = '  -1.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('E46').Value = '  -3.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '340.03'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.46'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.845'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.42%  '
